$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 4 and 5 hold two dataset entries (Populacao/Geracao combinations)
# that were saved out of sequence. Swap the two rows so the datasets end up
# with their columns "seguidas" (in sequence) with their neighbours.
# Columns C:G already hold identical values on both rows, so only A, B
# (text labels), I (numeric Average) and J (full-precision numeric-as-text
# StdDev) need to be exchanged between row 4 and row 5.
#
# Copy / PasteSpecial (values only) is used instead of a plain Value/Value2
# round-trip so that J's high-precision numeric-looking text keeps its
# original text cell-type (it must stay text, not become a Double) without
# picking up a new, unwanted number-format style in the process.
$xlPasteValues = -4163
$scratch = $ws.Range("Z1")

foreach ($col in @("A", "B", "I", "J")) {
    $cellRow4 = $ws.Range("${col}4")
    $cellRow5 = $ws.Range("${col}5")

    $cellRow4.Copy() | Out-Null
    $scratch.PasteSpecial($xlPasteValues) | Out-Null

    $cellRow5.Copy() | Out-Null
    $cellRow4.PasteSpecial($xlPasteValues) | Out-Null

    $scratch.Copy() | Out-Null
    $cellRow5.PasteSpecial($xlPasteValues) | Out-Null
}

$scratch.Clear() | Out-Null
$excel.CutCopyMode = 0

# Restore the selection to match the saved view state.
$ws.Range("B8").Select() | Out-Null
